$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.10920524597168
$ws.Range("B1").Value = 4.319094181060791
$ws.Range("C1").Value = 1.985294342041016
$ws.Range("D1").Value = 1.487732172012329
$ws.Range("E1").Value = 1.311786413192749
